# Atualizando a planilha de riscos e backlog
#
# Content-level edits applied to the "Riscos" table on Planilha1:
#  - B4: fix typo "Ferramente Nova" -> "Ferramenta Nova"
#  - G3: append extra guidance about talking to professors
#  - G5: "GESTÃO CONSTANTE E ATUALIZAÇÃO CONSTANTE" -> "GERENCIAMENTO E ATUALIZAÇÃO CONSTANTE"
#  - G6: "TRABALHAR FERRAMENTAS" -> "TRABALHAR COM FERRAMENTAS"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("B4").Value = "Ferramenta Nova"

$ws.Range("G3").Value = "PREPARAÇÃO PARA LIDAR COM O CONSTANTE APRENDIZADO QUE O PROJETO DEMANDA, COMUNICAÇÃO ENTRE A EQUIPE PARA SEGMENTAR TAREFAS E TRABALHO EM EQUIPE, ALÉM DE CONVERSAR COM OS PROFESSORES CASO NECESSÁRIO"

$ws.Range("G5").Value = "GERENCIAMENTO E ATUALIZAÇÃO CONSTANTE DO MODELO DE GESTÃO, PARA QUE AS FALHAS DE PLANEJAMENTO SEJAM NOTADAS O MAIS CEDO POSSÍVEL"

$ws.Range("G6").Value = "TRABALHAR COM FERRAMENTAS DE COMUNICAÇÃO EFETIVA ENTRE OS MEMBROS"
